$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Rows.Item(5).Delete()
